{"js": "// Apply the Korean wording tweaks to the \"Challenges\" bullet list paragraphs.\nconst replacements = [\n  {\n    find: \"\uc81c\ud55c\ub41c \ube0c\ub79c\ub4dc \uc778\uc2dd \ubc0f \uc778\uc2dd\",\n    replace: \"\uc81c\ud55c\ub41c \ube0c\ub79c\ub4dc \uc778\uc9c0\ub3c4 \ubc0f \uc778\uc2dd\"\n  },\n  {\n    find: \": \uc774\ub7ec\ud55c \uc0c8\ub85c\uc6b4 \uc2dc\uc7a5\uc5d0\uc11c \uac00\uc2dc\uc131\uc744 \ub2ec\uc131\ud558\ub294 \uac83\uc774 \uc8fc\uc694 \uc7a5\uc560\ubb3c\uc774\uba70, \ucc98\uc74c\ubd80\ud130 Adatum\uc758 \ube0c\ub79c\ub4dc \uc785\uc9c0\ub97c \uad6c\ucd95\ud558\uae30 \uc704\ud55c \uac15\ub825\ud55c \ub9c8\ucf00\ud305 \ub178\ub825\uc774 \ud544\uc694\ud569\ub2c8\ub2e4.\",\n    replace: \": \uc774\ub7ec\ud55c \uc0c8\ub85c\uc6b4 \uc2dc\uc7a5\uc5d0\uc11c \uac00\uc2dc\uc131\uc744 \ub2ec\uc131\ud558\ub294 \uac83\uc774 \uc8fc\uc694 \uacfc\uc81c\uc774\uba70, \ucc98\uc74c\ubd80\ud130 Adatum\uc758 \ube0c\ub79c\ub4dc \uc785\uc9c0\ub97c \uad6c\ucd95\ud558\uae30 \uc704\ud55c \uac15\ub825\ud55c \ub9c8\ucf00\ud305 \ub178\ub825\uc774 \ud544\uc694\ud569\ub2c8\ub2e4.\"\n  },\n  {\n    find: \": \uc774\ub7ec\ud55c \uc2dc\uc7a5\uc758 \ub2e4\uc591\ud55c \uc694\uad6c\uc5d0 \ub9de\uac8c \uc81c\ud488 \ubc0f \ub9c8\ucf00\ud305\uc744 \uc870\uc815\ud558\ub294 \uac83\uc740 \uc9c0\uc5ed \uae30\uc5c5 \ubc0f \uc18c\ube44\uc790\uc758 \uacf5\uac10\uc744 \uc774\ub04c\uc5b4 \ub0b4\ub294 \ub370 \ub9e4\uc6b0 \uc911\uc694\ud569\ub2c8\ub2e4.\",\n    replace: \": \uc774\ub7ec\ud55c \uc2dc\uc7a5\uc758 \ub2e4\uc591\ud55c \uc694\uad6c\uc5d0 \ub9de\uac8c \uc81c\ud488 \ubc0f \ub9c8\ucf00\ud305\uc744 \uc870\uc815\ud558\ub294 \uac83\uc740 \uc9c0\uc5ed \uae30\uc5c5 \ubc0f \uc18c\ube44\uc790\uc758 \uacf5\uac10\uc744 \uc774\ub04c\uc5b4 \ub0b4\ub294 \ub370 \uc788\uc5b4 \ub9e4\uc6b0 \uc911\uc694\ud569\ub2c8\ub2e4.\"\n  },\n  {\n    find: \": Adatum\uc740 \uc9c0\uc5ed\uc758 \uace0\uc720\ud55c \ub370\uc774\ud130 \uac1c\uc778 \uc815\ubcf4 \ubcf4\ud638, \ubcf4\uc548 \ubc0f \uc6b4\uc601 \uaddc\uc815\uc744 \ud0d0\uc0c9\ud558\ub294 \ubcf5\uc7a1\ud55c \uc791\uc5c5\uc5d0 \uc9c1\uba74\ud558\uace0 \uc788\uc73c\uba70, \ubd80\uc9c0\ub7f0\ud55c \uaddc\uc815 \uc900\uc218 \ub178\ub825\uc774 \ud544\uc694\ud569\ub2c8\ub2e4.\",\n    replace: \": Adatum\uc740 \uc9c0\uc5ed\uc758 \uace0\uc720\ud55c \ub370\uc774\ud130 \uac1c\uc778 \uc815\ubcf4 \ubcf4\ud638, \ubcf4\uc548 \ubc0f \uc6b4\uc601 \uaddc\uc815\uc744 \ud0d0\uc0c9\ud574\uc57c \ud558\ub294 \ubcf5\uc7a1\ud55c \uc791\uc5c5\uc5d0 \uc9c1\uba74\ud558\uace0 \uc788\uc73c\uba70, \ubd80\uc9c0\ub7f0\ud55c \uaddc\uc815 \uc900\uc218 \ub178\ub825\uc774 \ud544\uc694\ud569\ub2c8\ub2e4.\"\n  },\n  {\n    find: \": \ud6a8\uc728\uc801\uc778 \uc9c0\uc5ed \uac04 \uc6b4\uc601\uc744 \uc218\ub9bd\ud558\ub294 \uac83\uc740 \ud2b9\ud788 \ub192\uc740 \uc11c\ube44\uc2a4 \uc218\uc900\uc744 \uc720\uc9c0\ud558\uace0 \uc9c0\ub9ac\uc801 \uc704\uce58\uc5d0 \uac78\uccd0 \ub370\uc774\ud130 \uc13c\ud130\ub97c \uad00\ub9ac\ud558\ub294 \ub370 \ubb3c\ub958 \ubb38\uc81c\ub97c \uc81c\uc2dc\ud569\ub2c8\ub2e4.\",\n    replace: \": \ud6a8\uc728\uc801\uc778 \uc9c0\uc5ed \uac04 \uc6b4\uc601\uc744 \uc218\ub9bd\ud558\ub294 \uacfc\uc815\uc5d0\uc11c \ud2b9\ud788 \ub192\uc740 \uc11c\ube44\uc2a4 \uc218\uc900\uc744 \uc720\uc9c0\ud558\uace0 \uc9c0\ub9ac\uc801 \uc704\uce58\uc5d0 \uac78\uccd0 \ub370\uc774\ud130 \uc13c\ud130\ub97c \uad00\ub9ac\ud558\uae30 \uc704\ud574\uc11c\ub294 \ubb3c\ub958\uc0c1\uc758 \uc5b4\ub824\uc6c0\uc774 \uc788\uae30 \ub9c8\ub828\uc785\ub2c8\ub2e4.\"\n  }\n];\n\nfor (const { find, replace } of replacements) {\n  const results = context.document.body.search(find, { matchCase: true });\n  results.load(\"text\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text: \" + find);\n  }\n\n  for (const item of results.items) {\n    item.insertText(replace, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the Korean wording tweaks to the \"Challenges\" bullet list paragraphs.\n$d = $word.ActiveDocument\n\nfunction Replace-DocText($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $result = $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n    if (-not $result) {\n        throw \"Could not find text: $findText\"\n    }\n}\n\nReplace-DocText \"\uc81c\ud55c\ub41c \ube0c\ub79c\ub4dc \uc778\uc2dd \ubc0f \uc778\uc2dd\" \"\uc81c\ud55c\ub41c \ube0c\ub79c\ub4dc \uc778\uc9c0\ub3c4 \ubc0f \uc778\uc2dd\"\nReplace-DocText \": \uc774\ub7ec\ud55c \uc0c8\ub85c\uc6b4 \uc2dc\uc7a5\uc5d0\uc11c \uac00\uc2dc\uc131\uc744 \ub2ec\uc131\ud558\ub294 \uac83\uc774 \uc8fc\uc694 \uc7a5\uc560\ubb3c\uc774\uba70, \ucc98\uc74c\ubd80\ud130 Adatum\uc758 \ube0c\ub79c\ub4dc \uc785\uc9c0\ub97c \uad6c\ucd95\ud558\uae30 \uc704\ud55c \uac15\ub825\ud55c \ub9c8\ucf00\ud305 \ub178\ub825\uc774 \ud544\uc694\ud569\ub2c8\ub2e4.\" \": \uc774\ub7ec\ud55c \uc0c8\ub85c\uc6b4 \uc2dc\uc7a5\uc5d0\uc11c \uac00\uc2dc\uc131\uc744 \ub2ec\uc131\ud558\ub294 \uac83\uc774 \uc8fc\uc694 \uacfc\uc81c\uc774\uba70, \ucc98\uc74c\ubd80\ud130 Adatum\uc758 \ube0c\ub79c\ub4dc \uc785\uc9c0\ub97c \uad6c\ucd95\ud558\uae30 \uc704\ud55c \uac15\ub825\ud55c \ub9c8\ucf00\ud305 \ub178\ub825\uc774 \ud544\uc694\ud569\ub2c8\ub2e4.\"\nReplace-DocText \": \uc774\ub7ec\ud55c \uc2dc\uc7a5\uc758 \ub2e4\uc591\ud55c \uc694\uad6c\uc5d0 \ub9de\uac8c \uc81c\ud488 \ubc0f \ub9c8\ucf00\ud305\uc744 \uc870\uc815\ud558\ub294 \uac83\uc740 \uc9c0\uc5ed \uae30\uc5c5 \ubc0f \uc18c\ube44\uc790\uc758 \uacf5\uac10\uc744 \uc774\ub04c\uc5b4 \ub0b4\ub294 \ub370 \ub9e4\uc6b0 \uc911\uc694\ud569\ub2c8\ub2e4.\" \": \uc774\ub7ec\ud55c \uc2dc\uc7a5\uc758 \ub2e4\uc591\ud55c \uc694\uad6c\uc5d0 \ub9de\uac8c \uc81c\ud488 \ubc0f \ub9c8\ucf00\ud305\uc744 \uc870\uc815\ud558\ub294 \uac83\uc740 \uc9c0\uc5ed \uae30\uc5c5 \ubc0f \uc18c\ube44\uc790\uc758 \uacf5\uac10\uc744 \uc774\ub04c\uc5b4 \ub0b4\ub294 \ub370 \uc788\uc5b4 \ub9e4\uc6b0 \uc911\uc694\ud569\ub2c8\ub2e4.\"\nReplace-DocText \": Adatum\uc740 \uc9c0\uc5ed\uc758 \uace0\uc720\ud55c \ub370\uc774\ud130 \uac1c\uc778 \uc815\ubcf4 \ubcf4\ud638, \ubcf4\uc548 \ubc0f \uc6b4\uc601 \uaddc\uc815\uc744 \ud0d0\uc0c9\ud558\ub294 \ubcf5\uc7a1\ud55c \uc791\uc5c5\uc5d0 \uc9c1\uba74\ud558\uace0 \uc788\uc73c\uba70, \ubd80\uc9c0\ub7f0\ud55c \uaddc\uc815 \uc900\uc218 \ub178\ub825\uc774 \ud544\uc694\ud569\ub2c8\ub2e4.\" \": Adatum\uc740 \uc9c0\uc5ed\uc758 \uace0\uc720\ud55c \ub370\uc774\ud130 \uac1c\uc778 \uc815\ubcf4 \ubcf4\ud638, \ubcf4\uc548 \ubc0f \uc6b4\uc601 \uaddc\uc815\uc744 \ud0d0\uc0c9\ud574\uc57c \ud558\ub294 \ubcf5\uc7a1\ud55c \uc791\uc5c5\uc5d0 \uc9c1\uba74\ud558\uace0 \uc788\uc73c\uba70, \ubd80\uc9c0\ub7f0\ud55c \uaddc\uc815 \uc900\uc218 \ub178\ub825\uc774 \ud544\uc694\ud569\ub2c8\ub2e4.\"\nReplace-DocText \": \ud6a8\uc728\uc801\uc778 \uc9c0\uc5ed \uac04 \uc6b4\uc601\uc744 \uc218\ub9bd\ud558\ub294 \uac83\uc740 \ud2b9\ud788 \ub192\uc740 \uc11c\ube44\uc2a4 \uc218\uc900\uc744 \uc720\uc9c0\ud558\uace0 \uc9c0\ub9ac\uc801 \uc704\uce58\uc5d0 \uac78\uccd0 \ub370\uc774\ud130 \uc13c\ud130\ub97c \uad00\ub9ac\ud558\ub294 \ub370 \ubb3c\ub958 \ubb38\uc81c\ub97c \uc81c\uc2dc\ud569\ub2c8\ub2e4.\" \": \ud6a8\uc728\uc801\uc778 \uc9c0\uc5ed \uac04 \uc6b4\uc601\uc744 \uc218\ub9bd\ud558\ub294 \uacfc\uc815\uc5d0\uc11c \ud2b9\ud788 \ub192\uc740 \uc11c\ube44\uc2a4 \uc218\uc900\uc744 \uc720\uc9c0\ud558\uace0 \uc9c0\ub9ac\uc801 \uc704\uce58\uc5d0 \uac78\uccd0 \ub370\uc774\ud130 \uc13c\ud130\ub97c \uad00\ub9ac\ud558\uae30 \uc704\ud574\uc11c\ub294 \ubb3c\ub958\uc0c1\uc758 \uc5b4\ub824\uc6c0\uc774 \uc788\uae30 \ub9c8\ub828\uc785\ub2c8\ub2e4.\"\n"}
